# Working on new CategoryView - Scene switch only
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Barry Greene's URL moved to his new video-lesson platform
$ws.Range("C6").Value = "https://barrygreenevideolessons.vhx.tv/"

# New row: JavaFX mobile-apps article (Gluon / foojay.io)
# "15" must stay a text value (like the rest of column A), so force it
# with a leading apostrophe and then drop the resulting text-format style
# so the cell matches the unstyled look of the other data rows.
$ws.Range("A13").Value = "'15"
$ws.Range("A13").ClearFormats()
$ws.Range("B13").Value = "Creating mobile apps with JavaFX"
$ws.Range("C13").Value = "https://foojay.io/today/creating-mobile-apps-with-javafx-part-1/"
$ws.Range("D13").Value = "How to run a JavaFX programme on mobiles using Gluon"
$ws.Range("E13").Value = "Coding"

# Column B widened slightly to fit the new longest title
$ws.Columns.Item(2).ColumnWidth = 30
